$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.385.62'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.66%  '
$ws.Range("E2").NumberFormat = "General"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.653.87'
$ws.Range("D3").NumberFormat = "General"

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("E4").NumberFormat = "General"

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '213.44'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.45%  '
$ws.Range("E5").NumberFormat = "General"

# Row 6
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("E6").NumberFormat = "General"

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").NumberFormat = "General"

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("E8").NumberFormat = "General"

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.261'
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.17%  '
$ws.Range("E9").NumberFormat = "General"

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0615'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.82%  '
$ws.Range("E10").NumberFormat = "General"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0876'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.30%  '
$ws.Range("E11").NumberFormat = "General"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.887.60'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.56%  '
$ws.Range("E12").NumberFormat = "General"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.647.79'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -1.05%  '
$ws.Range("E13").NumberFormat = "General"

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.573'
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +4.28%  '
$ws.Range("E14").NumberFormat = "General"

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.07'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("E15").NumberFormat = "General"

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.68'
$ws.Range("D16").NumberFormat = "General"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '27.365.88'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.62%  '
$ws.Range("E17").NumberFormat = "General"

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '232.10'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -6.72%  '
$ws.Range("E18").NumberFormat = "General"

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0₃0726'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.68%  '
$ws.Range("E19").NumberFormat = "General"

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.52'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.10%  '
$ws.Range("E20").NumberFormat = "General"

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.39%  '
$ws.Range("E21").NumberFormat = "General"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.95%  '
$ws.Range("E22").NumberFormat = "General"

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.21'
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -1.48%  '
$ws.Range("E23").NumberFormat = "General"

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.11%  '
$ws.Range("E24").NumberFormat = "General"

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '146.97'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.22%  '
$ws.Range("E25").NumberFormat = "General"

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.18'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.42%  '
$ws.Range("E26").NumberFormat = "General"

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.89'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.82%  '
$ws.Range("E27").NumberFormat = "General"

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.34%  '
$ws.Range("E28").NumberFormat = "General"

# Row 29
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.02%  '
$ws.Range("E29").NumberFormat = "General"

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0498'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.12%  '
$ws.Range("E30").NumberFormat = "General"

# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -3.54%  '
$ws.Range("E31").NumberFormat = "General"

# Row 32
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.87%  '
$ws.Range("E32").NumberFormat = "General"

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.458.34'
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.13%  '
$ws.Range("E33").NumberFormat = "General"

# Row 34
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.52%  '
$ws.Range("E34").NumberFormat = "General"

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("E35").NumberFormat = "General"

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.38'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("E36").NumberFormat = "General"

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.29%  '
$ws.Range("E37").NumberFormat = "General"

# Row 38
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.10%  '
$ws.Range("E38").NumberFormat = "General"

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.53%  '
$ws.Range("E40").NumberFormat = "General"

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.29%  '
$ws.Range("E41").NumberFormat = "General"

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.47'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.93%  '
$ws.Range("E42").NumberFormat = "General"

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '65.16'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -5.77%  '
$ws.Range("E43").NumberFormat = "General"

# Row 44
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.36%  '
$ws.Range("E44").NumberFormat = "General"

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.795.87'
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.45%  '
$ws.Range("E45").NumberFormat = "General"

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.786'
$ws.Range("D46").NumberFormat = "General"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.56%  '
$ws.Range("E46").NumberFormat = "General"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.05%  '
$ws.Range("E47").NumberFormat = "General"

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '88.12'
$ws.Range("D48").NumberFormat = "General"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -0.26%  '
$ws.Range("E48").NumberFormat = "General"

# Row 49
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.17%  '
$ws.Range("E49").NumberFormat = "General"

# Row 50
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("B50").NumberFormat = "General"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C50").NumberFormat = "General"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.75'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.23%  '
$ws.Range("E50").NumberFormat = "General"

# Row 51
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("B51").NumberFormat = "General"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("C51").NumberFormat = "General"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₇0984'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -7.46%  '
$ws.Range("E51").NumberFormat = "General"
